$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.302.17'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.10%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.618.89'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.74%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.99'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.57%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.484'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.60%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.44%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.54%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.78'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +4.48%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0815'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.83%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.76%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.621.72'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.88%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.02'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.90%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.518'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.23%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.297.49'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.05%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.38'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +3.66%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.51%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '201.81'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.33'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.41%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.51%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -3.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.45'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.46%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.25%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.17'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.32%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.32%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0514'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +8.39%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.62%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.92%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.50'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.55%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +2.49%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.177.11'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +4.26%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.53%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.809'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.40%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.02%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.04%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.24%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.36'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +4.35%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.39%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.757.22'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.92%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.46%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.65%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '53.75'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.32%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0507'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.92%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.53%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.22%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.13%  '
